# Adds a "count" column (letter frequency counts) to Sheet1 and Sheet2,
# plus a summary/labelled-count block and grand total on Sheet2.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$counts = @(141, 23, 102, 8, 60, 9, 41, 8, 4, 37, 4, 131, 74, 1, 82, 84, 124, 71, 6, 70)

# --- Seed the new shared strings in the same order as the target file:
#     "letter" first, then "count". ---
$ws2.Range("E1").Value = "letter"
$ws1.Range("C1").Value = "count"

# --- Sheet1: new column C ("count") ---
for ($i = 0; $i -lt $counts.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 3).Value = $counts[$i]
}

$ws1.Range("C21").Select()

# --- Sheet2: new columns E ("letter"), F ("count"), G (concat formula) + grand total ---
$ws2.Range("F1").Value = "count"

for ($i = 0; $i -lt $counts.Length; $i++) {
    $row = $i + 2
    $ws2.Range("E$row").Formula = "=A$row"
    $ws2.Cells.Item($row, 6).Value = $counts[$i]
    $ws2.Range("G$row").Formula = '=CONCATENATE(A' + $row + ',"     ",F' + $row + ')'
}

$ws2.Range("F22").Formula = "=SUM(F2:F21)"

$ws2.Range("F2:F21").Select()
